$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.530.68"
$ws.Range("E2").Value = "  +2.88%  "
$ws.Range("D3").Value = "1.603.40"
$ws.Range("E3").Value = "  +2.28%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("E6").Value = "  +6.87%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "26.82"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.50"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.31%  "
$ws.Range("E10").Value = "  +2.50%  "
$ws.Range("E11").Value = "  +2.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0909"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.62%  "
$ws.Range("D13").Value = "1.832.99"
$ws.Range("E13").Value = "  +2.35%  "
$ws.Range("D14").Value = "1.590.07"
$ws.Range("E14").Value = "  +1.49%  "
$ws.Range("D15").Value = "29.511.85"
$ws.Range("E15").Value = "  +2.83%  "
$ws.Range("E16").Value = "  +3.48%  "
$ws.Range("E17").Value = "  +1.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.26%  "
$ws.Range("D21").Value = "0.0₃0691"
$ws.Range("E21").Value = "  +1.55%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.72%  "
$ws.Range("E27").Value = "  +5.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.06%  "
$ws.Range("E29").Value = "  +2.37%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("E31").Value = "  +2.61%  "
$ws.Range("E32").Value = "  +0.82%  "
$ws.Range("E33").Value = "  +1.49%  "
$ws.Range("E34").Value = "  +3.33%  "
$ws.Range("D35").Value = "1.409.27"
$ws.Range("E35").Value = "  +1.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.59%  "
$ws.Range("E37").Value = "  +3.87%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.81"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.28%  "
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("E40").Value = "  +2.23%  "
$ws.Range("E41").Value = "  +3.67%  "
$ws.Range("E42").Value = "  +0.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0489"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "53.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +23.10%  "
$ws.Range("E45").Value = "  +3.04%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "65.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.28"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.71%  "
$ws.Range("D49").Value = "1.742.32"
$ws.Range("E49").Value = "  +2.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.856"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "86.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.44%  "
